$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing data rows (31..112) down to (32..113)
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with a new weekly record (same static attributes as the
# rest of the "Jengibre" series, with its own date / volume / price values)
$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44998
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 100114007
$ws.Range("G31").Value = "Jengibre"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 440
$ws.Range("K31").Value = 17500
$ws.Range("L31").Value = 18000
$ws.Range("M31").Value = 17750
$ws.Range("N31").Value = "`$/caja 13 kilos"
$ws.Range("O31").Value = "Perú"
$ws.Range("P31").Value = 1365
$ws.Range("Q31").Value = 13
$ws.Range("R31").Value = "Hortaliza"
